# Implement foreign riders to riders list:
# append one new rider (Peter PALASTHY, SVK) as the new last row (768)
# of the riders list sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 768

# Licence_num / UCI_ID / FederationID / International Licence Code
$ws.Cells.Item($row, 1).Value = 10005129361
$ws.Cells.Item($row, 2).Value = 10005129361
$ws.Cells.Item($row, 4).Value = 10005129361
$ws.Cells.Item($row, 5).Value = 10005129361

# Expiry_date - force Text format first so the date-looking string isn't
# auto-converted into a date serial by Excel's input parsing.
$ws.Cells.Item($row, 6).NumberFormat = "@"
$ws.Cells.Item($row, 6).Value = "2022/12/31"

# Licence_type
$ws.Cells.Item($row, 7).Value = "BMX-RACE"

# Dob - same text-format trick as Expiry_date above.
$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 9).Value = "1981/02/20"

# First_name / Surname
$ws.Cells.Item($row, 10).Value = "Peter"
$ws.Cells.Item($row, 11).Value = "PALÁSTHY"

# Sex / CLUB / State / UCI_Country / Class
$ws.Cells.Item($row, 16).Value = "M"
$ws.Cells.Item($row, 17).Value = "Slovakia-All Clubs"
$ws.Cells.Item($row, 18).Value = "SVK"
$ws.Cells.Item($row, 19).Value = "CZE"
$ws.Cells.Item($row, 20).Value = "Boys 11/12"

# Plate / Plate2
$ws.Cells.Item($row, 24).Value = 0
$ws.Cells.Item($row, 25).Value = 0

# Ranking / Ranking2 - present but blank (text) cells, same as other rows.
$ws.Cells.Item($row, 28).NumberFormat = "@"
$ws.Cells.Item($row, 28).Value = ""
$ws.Cells.Item($row, 29).NumberFormat = "@"
$ws.Cells.Item($row, 29).Value = ""

# columns AJ / AK template markers
$ws.Cells.Item($row, 36).Value = "T1"
$ws.Cells.Item($row, 37).Value = "T2"
